$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.275.49"
$ws.Range("E2").Value = "  +0.62%  "

$ws.Range("D3").Value = "1.666.09"
$ws.Range("E3").Value = "  +0.70%  "

$ws.Range("E4").Value = "  +0.76%  "

$ws.Range("E5").Value = "  +0.47%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5327"
$ws.Range("D6").Style = $ws.Range("B6").Style
$ws.Range("E6").Value = "  +1.76%  "

$ws.Range("E7").Value = "  +0.73%  "

$ws.Range("E8").Value = "  +1.30%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06381"
$ws.Range("D9").Style = $ws.Range("B9").Style
$ws.Range("E9").Value = "  +0.59%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.54"
$ws.Range("D10").Style = $ws.Range("B10").Style
$ws.Range("E10").Value = "  +0.96%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07830"
$ws.Range("D11").Style = $ws.Range("B11").Style
$ws.Range("E11").Value = "  +0.41%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.560"
$ws.Range("D12").Style = $ws.Range("B12").Style
$ws.Range("E12").Value = "  +1.38%  "

$ws.Range("D13").Value = "1.664.69"
$ws.Range("E13").Value = "  +1.22%  "

$ws.Range("D14").Value = "1.895.03"
$ws.Range("E14").Value = "  +0.70%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5531"
$ws.Range("D15").Style = $ws.Range("B15").Style
$ws.Range("E15").Value = "  +1.16%  "

$ws.Range("D16").Value = "0.0₅8213"
$ws.Range("E16").Value = "  +0.24%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.74"
$ws.Range("D17").Style = $ws.Range("B17").Style
$ws.Range("E17").Value = "  +0.67%  "

$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "26.284.77"
$ws.Range("E18").Value = "  +0.65%  "

$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.010"
$ws.Range("D19").Style = $ws.Range("B19").Style
$ws.Range("E19").Value = "  +0.78%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.678"
$ws.Range("D20").Style = $ws.Range("B20").Style
$ws.Range("E20").Value = "  +2.26%  "

$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.48"
$ws.Range("D21").Style = $ws.Range("B21").Style
$ws.Range("E21").Value = "  +1.15%  "

$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.22"
$ws.Range("D22").Style = $ws.Range("B22").Style
$ws.Range("E22").Value = "  +1.62%  "

$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.033"
$ws.Range("D23").Style = $ws.Range("B23").Style
$ws.Range("E23").Value = "  +0.09%  "

$ws.Range("B24").Value = "BinanceUSD"
$ws.Range("C24").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.011"
$ws.Range("D24").Style = $ws.Range("B24").Style
$ws.Range("E24").Value = "  +0.70%  "

$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.64"
$ws.Range("D25").Style = $ws.Range("B25").Style
$ws.Range("E25").Value = "  +2.51%  "

$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1225"
$ws.Range("D26").Style = $ws.Range("B26").Style
$ws.Range("E26").Value = "  -1.14%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.200"
$ws.Range("D27").Style = $ws.Range("B27").Style
$ws.Range("E27").Value = "  -0.47%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.12"
$ws.Range("D28").Style = $ws.Range("B28").Style
$ws.Range("E28").Value = "  -0.07%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.483"
$ws.Range("D29").Style = $ws.Range("B29").Style
$ws.Range("E29").Value = "  +3.83%  "

$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05895"
$ws.Range("D30").Style = $ws.Range("B30").Style
$ws.Range("E30").Value = "  -0.07%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.279"
$ws.Range("D31").Style = $ws.Range("B31").Style
$ws.Range("E31").Value = "  +0.17%  "

$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.601"
$ws.Range("D32").Style = $ws.Range("B32").Style
$ws.Range("E32").Value = "  +2.27%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.277"
$ws.Range("D33").Style = $ws.Range("B33").Style
$ws.Range("E33").Value = "  +0.99%  "

$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.613"
$ws.Range("D34").Style = $ws.Range("B34").Style
$ws.Range("E34").Value = "  +1.58%  "

$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9637"
$ws.Range("D35").Style = $ws.Range("B35").Style
$ws.Range("E35").Value = "  +1.21%  "

$ws.Range("B36").Value = "MXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.824"
$ws.Range("D36").Style = $ws.Range("B36").Style
$ws.Range("E36").Value = "  +1.35%  "

$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.423"
$ws.Range("D37").Style = $ws.Range("B37").Style
$ws.Range("E37").Value = "  +0.56%  "

$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5798"
$ws.Range("D38").Style = $ws.Range("B38").Style
$ws.Range("E38").Value = "  +2.07%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01606"
$ws.Range("D39").Style = $ws.Range("B39").Style
$ws.Range("E39").Value = "  -0.64%  "

$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8632"
$ws.Range("D40").Style = $ws.Range("B40").Style
$ws.Range("E40").Value = "  +1.75%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.840"
$ws.Range("D41").Style = $ws.Range("B41").Style
$ws.Range("E41").Value = "  +0.28%  "

$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.009"
$ws.Range("D42").Style = $ws.Range("B42").Style
$ws.Range("E42").Value = "  +0.67%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.047.86"
$ws.Range("E43").Value = "  +1.58%  "

$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "103.99"
$ws.Range("D44").Style = $ws.Range("B44").Style
$ws.Range("E44").Value = "  +1.05%  "

$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "1.805.53"
$ws.Range("E45").Value = "  +0.43%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.72"
$ws.Range("D46").Style = $ws.Range("B46").Style
$ws.Range("E46").Value = "  +1.02%  "

$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.009"
$ws.Range("D47").Style = $ws.Range("B47").Style
$ws.Range("E47").Value = "  +0.42%  "

$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₈105"
$ws.Range("E48").Value = "  -5.99%  "

$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4380"
$ws.Range("D49").Style = $ws.Range("B49").Style
$ws.Range("E49").Value = "  +1.69%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.037"
$ws.Range("D50").Style = $ws.Range("B50").Style
$ws.Range("E50").Value = "  +2.49%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05161"
$ws.Range("D51").Style = $ws.Range("B51").Style
$ws.Range("E51").Value = "  -0.06%  "
